# Apply the edits described by the diff to the "Metadata" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental"): B7 was blank -> "false".
# A literal Range.Value assignment of the text "false" gets auto-typed as a
# Boolean by Excel's smart-input parsing (t="b"), but the source file (a FHIR
# IG publisher export) always stores plain shared-string text (t="s"). Build
# the text as a formula result first, then paste back as a value-only copy so
# the literal string "false" lands in the cell without being reinterpreted as
# a boolean.
$cell = $ws.Cells.Item(7, 2)
$cell.Formula = "=""fal""&""se"""
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Row 8 ("Date"): updated timestamp.
$ws.Cells.Item(8, 2).Value = "2025-11-30T13:08:37+00:00"

# Row 17 ("Description"): was blank -> descriptive text.
$ws.Cells.Item(17, 2).Value = "Codes for recommended activity intensity based on recovery status"
